# Update gh-pages to output generated at 456a3b4
# Applies updated "想去人数" (F column) counts across the four sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 34
$ws1.Range("F5").Value  = 959
$ws1.Range("F6").Value  = 356
$ws1.Range("F8").Value  = 553
$ws1.Range("F9").Value  = 1434
$ws1.Range("F11").Value = 1326
$ws1.Range("F13").Value = 387
$ws1.Range("F14").Value = 1602
$ws1.Range("F15").Value = 1357
$ws1.Range("F16").Value = 782
$ws1.Range("F17").Value = 232
$ws1.Range("F18").Value = 1362
$ws1.Range("F19").Value = 261
$ws1.Range("F20").Value = 63
$ws1.Range("F21").Value = 1113
$ws1.Range("F22").Value = 394
$ws1.Range("F23").Value = 3451
$ws1.Range("F24").Value = 674
$ws1.Range("F25").Value = 560
$ws1.Range("F26").Value = 1527

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F8").Value  = 19
$ws2.Range("F12").Value = 76

# --- Sheet "本地生活" ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 793

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 34
$ws4.Range("F3").Value  = 793
$ws4.Range("F13").Value = 19
$ws4.Range("F15").Value = 959
$ws4.Range("F16").Value = 356
$ws4.Range("F18").Value = 553
$ws4.Range("F19").Value = 1434
$ws4.Range("F21").Value = 1326
$ws4.Range("F23").Value = 387
$ws4.Range("F24").Value = 1602
$ws4.Range("F25").Value = 1357
$ws4.Range("F26").Value = 782
$ws4.Range("F27").Value = 232
$ws4.Range("F28").Value = 1362
$ws4.Range("F29").Value = 261
$ws4.Range("F30").Value = 63
$ws4.Range("F33").Value = 1113
$ws4.Range("F34").Value = 394
$ws4.Range("F35").Value = 3451
$ws4.Range("F36").Value = 674
$ws4.Range("F37").Value = 560
$ws4.Range("F38").Value = 1527
$ws4.Range("F39").Value = 76
